$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.106.24'
$ws.Range("E2").Value = '  -3.74%  '
$ws.Range("D3").Value = '2.460.68'
$ws.Range("E3").Value = '  -3.00%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''311.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '''94.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.38%  '
$ws.Range("D7").Value = '''0.550'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.79%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").Value = '''0.502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.23%  '
$ws.Range("D10").Value = '''33.44'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.97%  '
$ws.Range("E11").Value = '  -3.07%  '
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '''6.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.62%  '
$ws.Range("D14").Value = '2.831.92'
$ws.Range("E14").Value = '  -3.41%  '
$ws.Range("D15").Value = '2.453.55'
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '''14.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.67%  '
$ws.Range("D17").Value = '''0.787'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.93%  '
$ws.Range("D18").Value = '41.070.13'
$ws.Range("E18").Value = '  -3.83%  '
$ws.Range("D19").Value = '''6.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.72%  '
$ws.Range("D20").Value = '0.0₃0916'
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").Value = '''11.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.98%  '
$ws.Range("D22").Value = '''67.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.54%  '
$ws.Range("D23").Value = '''236.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.48%  '
$ws.Range("E24").Value = '  -4.37%  '
$ws.Range("D25").Value = '''1.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.01%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").Value = '''24.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.22%  '
$ws.Range("E28").Value = '  -5.12%  '
$ws.Range("D29").Value = '''9.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.04%  '
$ws.Range("E30").Value = '  -7.66%  '
$ws.Range("D31").Value = '''152.95'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.49%  '
$ws.Range("D32").Value = '''5.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.73%  '
$ws.Range("D33").Value = '''2.58'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("D34").Value = '''2.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.12%  '
$ws.Range("D35").Value = '''0.0752'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.34%  '
$ws.Range("D36").Value = '''3.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.75%  '
$ws.Range("E37").Value = '  -6.73%  '
$ws.Range("D38").Value = '''17.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.05%  '
$ws.Range("E39").Value = '  -8.13%  '
$ws.Range("E40").Value = '  -4.22%  '
$ws.Range("D41").Value = '''4.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("D42").Value = '''21.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.30%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '1.966.08'
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("D46").Value = '''3.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.95%  '
$ws.Range("D47").Value = '''8.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").Value = '''69.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.97%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").Value = '''76.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.59%  '
$ws.Range("D50").Value = '''97.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("E51").Value = '  -6.55%  '
